# Insert a new weekly price record as row 498 (pushing the existing
# rows 498:527 down to 499:528) on the single data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 498..527 down to 499..528, leaving row 498 blank.
$ws.Rows.Item(498).Insert()

# Populate the newly inserted row 498 with the new weekly record.
$ws.Cells.Item(498, 1).Value  = 6
$ws.Cells.Item(498, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(498, 3).Value  = "Metropolitana"
$ws.Cells.Item(498, 4).Value  = 44747
$ws.Cells.Item(498, 5).Value  = 13
$ws.Cells.Item(498, 6).Value  = 100112044
$ws.Cells.Item(498, 7).Value  = "Perejil"
$ws.Cells.Item(498, 8).Value  = "Sin especificar"
$ws.Cells.Item(498, 9).Value  = "Primera"
$ws.Cells.Item(498, 10).Value = 140
$ws.Cells.Item(498, 11).Value = 17000
$ws.Cells.Item(498, 12).Value = 18000
$ws.Cells.Item(498, 13).Value = 17429
$ws.Cells.Item(498, 14).Value = "$/docena de atados"
$ws.Cells.Item(498, 15).Value = "Región Metropolitana"
$ws.Cells.Item(498, 16).Value = 5810
$ws.Cells.Item(498, 17).Value = 3
$ws.Cells.Item(498, 18).Value = "Hortaliza"
